# feat: add 2022-Q1 data
#
# - Insert a new "2022-Q1" worksheet (fund holdings detail) right before the
#   "总计" (summary) worksheet.
# - Add a new summary row for "2022-Q1" at the top of the "总计" worksheet,
#   pushing the previously existing rows down by one.

$wb = $excel.ActiveWorkbook

# Use an existing, similarly-laid-out sheet as a formatting template for the
# header row / index column styling (NOTE: always resolve "总计" freshly via
# Worksheets.Item(...) instead of caching it in a variable that also gets
# passed into Worksheets.Add(), since that rebinds the variable to the new
# sheet).
$templateSheet = $wb.Worksheets.Item("2021-Q4")

# --- Create the new "2022-Q1" worksheet right before "总计" ---
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$newSheet.Name = "2022-Q1"

# --- Header row (text) ---
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Match header styling used by the other quarter sheets.
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# --- Fund holdings data rows ---
# (Code / name / size / total position / position ratio / market value are
# stored as plain text in the source data, matching the rest of the workbook;
# only the row index (col A) and position rank (col H) are numeric. Force
# column B:G to Text format first so the numeric-looking strings are not
# silently reinterpreted as numbers.)
$newSheet.Range("B2:G7").NumberFormat = "@"

$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(2, 2).Value = "516150"
$newSheet.Cells.Item(2, 3).Value = "嘉实中证稀土产业ETF"
$newSheet.Cells.Item(2, 4).Value = "25.17"
$newSheet.Cells.Item(2, 5).Value = "99.75"
$newSheet.Cells.Item(2, 6).Value = "4.05"
$newSheet.Cells.Item(2, 7).Value = "1.0194"
$newSheet.Cells.Item(2, 8).Value = 8

$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(3, 2).Value = "516780"
$newSheet.Cells.Item(3, 3).Value = "华泰柏瑞中证稀土产业ETF"
$newSheet.Cells.Item(3, 4).Value = "11.06"
$newSheet.Cells.Item(3, 5).Value = "98.70"
$newSheet.Cells.Item(3, 6).Value = "4.01"
$newSheet.Cells.Item(3, 7).Value = "0.4435"
$newSheet.Cells.Item(3, 8).Value = 8

$newSheet.Cells.Item(4, 1).Value = 2
$newSheet.Cells.Item(4, 2).Value = "159715"
$newSheet.Cells.Item(4, 3).Value = "易方达中证稀土产业ETF"
$newSheet.Cells.Item(4, 4).Value = "3.42"
$newSheet.Cells.Item(4, 5).Value = "99.06"
$newSheet.Cells.Item(4, 6).Value = "4.00"
$newSheet.Cells.Item(4, 7).Value = "0.1368"
$newSheet.Cells.Item(4, 8).Value = 8

$newSheet.Cells.Item(5, 1).Value = 3
$newSheet.Cells.Item(5, 2).Value = "159713"
$newSheet.Cells.Item(5, 3).Value = "富国中证稀土产业交易型开放式指数证券投资基金"
$newSheet.Cells.Item(5, 4).Value = "3.26"
$newSheet.Cells.Item(5, 5).Value = "99.26"
$newSheet.Cells.Item(5, 6).Value = "4.03"
$newSheet.Cells.Item(5, 7).Value = "0.1314"
$newSheet.Cells.Item(5, 8).Value = 8

$newSheet.Cells.Item(6, 1).Value = 4
$newSheet.Cells.Item(6, 2).Value = "014331"
$newSheet.Cells.Item(6, 3).Value = "华泰柏瑞中证稀土产业ETF联接A"
$newSheet.Cells.Item(6, 4).Value = "0.86"
$newSheet.Cells.Item(6, 5).Value = "24.22"
$newSheet.Cells.Item(6, 6).Value = "1.08"
$newSheet.Cells.Item(6, 7).Value = "0.0093"
$newSheet.Cells.Item(6, 8).Value = 8

$newSheet.Cells.Item(7, 1).Value = 5
$newSheet.Cells.Item(7, 2).Value = "014332"
$newSheet.Cells.Item(7, 3).Value = "华泰柏瑞中证稀土产业ETF联接C"
$newSheet.Cells.Item(7, 4).Value = "0.70"
$newSheet.Cells.Item(7, 5).Value = "24.22"
$newSheet.Cells.Item(7, 6).Value = "1.08"
$newSheet.Cells.Item(7, 7).Value = "0.0076"
$newSheet.Cells.Item(7, 8).Value = 8

# The values are now stored as text (because of the "@" format applied
# above); strip the format/style again so the cells end up unstyled, same
# as in the other quarter sheets.
$newSheet.Range("B2:G7").ClearFormats()

# Apply the index-column styling (col A) used by the other quarter sheets.
$templateSheet.Range("A2").Copy()
$newSheet.Range("A2:A7").PasteSpecial(-4122)

# --- Update the "总计" (summary) sheet: add a 2022-Q1 row on top ---
$totalSheet = $wb.Worksheets.Item("总计")

# Push the 4 existing data rows down by one (writing literal values directly,
# bottom-up, so no row ends up temporarily holding another row's data).
$totalSheet.Cells.Item(6, 1).Value = 4
$totalSheet.Cells.Item(6, 2).Value = "2020-Q4"
$totalSheet.Cells.Item(6, 3).Value = 41
$totalSheet.Cells.Item(6, 4).Value = 14.03

$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(5, 2).Value = "2021-Q1"
$totalSheet.Cells.Item(5, 3).Value = 15
$totalSheet.Cells.Item(5, 4).Value = 2.4

$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(4, 2).Value = "2021-Q2"
$totalSheet.Cells.Item(4, 3).Value = 6
$totalSheet.Cells.Item(4, 4).Value = 2.8

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 2).Value = "2021-Q4"
$totalSheet.Cells.Item(3, 3).Value = 4
$totalSheet.Cells.Item(3, 4).Value = 2.01

# New row for the current quarter.
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 6
$totalSheet.Cells.Item(2, 4).Value = 1.75

# Row 6 is brand new territory on this sheet; give its index cell (col A)
# the same styling as the other index cells.
$totalSheet.Range("A5").Copy()
$totalSheet.Range("A6").PasteSpecial(-4122)
$totalSheet.Cells.Item(6, 1).Value = 4
